$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# forceMinRigidity / forceMaxRigidity: recomputed values
# ("all is computed right except non-linearity")
$ws.Range("B7").Value = 2.63
$ws.Range("B8").Value = 1.315

# Move the active selection to C10
$ws.Range("C10").Select()

# Reposition/resize the workbook window to match the saved view state
$win = $excel.ActiveWindow
$win.Left = 12285
$win.Top = 795
$win.Width = 16050
$win.Height = 13680
